# Update the "Forecast coal outages" row (row 5) on the "Coal outages" sheet
# so the later-day columns mirror the now-known values, and correct the
# Loy Yang A 4 outage return-date row (row 16) with its real expected
# return date / days-until-return.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coal outages")

# Row 5: forecast coal outages (in MW) per day
$ws.Range("E5").Value = 1465
$ws.Range("F5").Value = 1465
$ws.Range("G5").Value = 935
$ws.Range("H5").Value = 935

# Row 16: Loy Yang A 4 - expected return date (serial) and days until return
$ws.Range("G16").Value = 45642
$ws.Range("H16").Value = 11
